$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.196.31'
$ws.Range("E2").Value = '  +2.48%  '
$ws.Range("D3").Value = '3.265.76'
$ws.Range("E3").Value = '  +1.49%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '397.56'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.73'
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.582'
$ws.Range("E7").Value = '  +4.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.622'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.50'
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0955'
$ws.Range("E11").Value = '  +6.15%  '
$ws.Range("E12").Value = '  +1.90%  '
$ws.Range("D13").Value = '3.780.34'
$ws.Range("E13").Value = '  +1.46%  '
$ws.Range("E14").Value = '  +2.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.01'
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").Value = '3.273.56'
$ws.Range("E16").Value = '  +1.44%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.03'
$ws.Range("E17").Value = '  -1.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.13'
$ws.Range("E18").Value = '  +4.61%  '
$ws.Range("D19").Value = '57.000.72'
$ws.Range("E19").Value = '  +2.26%  '
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000109'
$ws.Range("E21").Value = '  +6.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.97'
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '292.74'
$ws.Range("E23").Value = '  -3.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.28'
$ws.Range("E24").Value = '  -1.39%  '
$ws.Range("E25").Value = '  -1.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.18'
$ws.Range("E26").Value = '  -0.23%  '
$ws.Range("E27").Value = '  -3.44%  '
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.40'
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("E30").Value = '  -2.90%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  +1.52%  '
$ws.Range("E33").Value = '  -1.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '39.92'
$ws.Range("E34").Value = '  +10.30%  '
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.47'
$ws.Range("E39").Value = '  -1.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.03'
$ws.Range("E40").Value = '  -2.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '137.00'
$ws.Range("E41").Value = '  +1.74%  '
$ws.Range("E42").Value = '  +1.81%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.94'
$ws.Range("E43").Value = '  -2.86%  '
$ws.Range("E44").Value = '  -2.59%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.283'
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.76'
$ws.Range("E46").Value = '  -1.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.47'
$ws.Range("E47").Value = '  +0.66%  '
$ws.Range("E48").Value = '  +4.66%  '
$ws.Range("D49").Value = '2.152.25'
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("E50").Value = '  -5.59%  '
$ws.Range("E51").Value = '  -5.81%  '
